$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = "42.896.13"
$r.Style = "Normal"

$r = $ws.Range("E2")
$r.NumberFormat = "@"
$r.Value = "  -1.66%  "
$r.Style = "Normal"

$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = "2.557.99"
$r.Style = "Normal"

$r = $ws.Range("E3")
$r.NumberFormat = "@"
$r.Value = "  -1.15%  "
$r.Style = "Normal"

$r = $ws.Range("E4")
$r.NumberFormat = "@"
$r.Value = "  +0.04%  "
$r.Style = "Normal"

$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "302.01"
$r.Style = "Normal"

$r = $ws.Range("E5")
$r.NumberFormat = "@"
$r.Value = "  +0.33%  "
$r.Style = "Normal"

$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "92.86"
$r.Style = "Normal"

$r = $ws.Range("E6")
$r.NumberFormat = "@"
$r.Value = "  -3.67%  "
$r.Style = "Normal"

$r = $ws.Range("E7")
$r.NumberFormat = "@"
$r.Value = "  -0.58%  "
$r.Style = "Normal"

$r = $ws.Range("E8")
$r.NumberFormat = "@"
$r.Value = "  -0.08%  "
$r.Style = "Normal"

$r = $ws.Range("E9")
$r.NumberFormat = "@"
$r.Value = "  -2.03%  "
$r.Style = "Normal"

$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "36.21"
$r.Style = "Normal"

$r = $ws.Range("E10")
$r.NumberFormat = "@"
$r.Value = "  -2.50%  "
$r.Style = "Normal"

$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "0.0811"
$r.Style = "Normal"

$r = $ws.Range("E11")
$r.NumberFormat = "@"
$r.Value = "  -0.44%  "
$r.Style = "Normal"

$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "7.78"
$r.Style = "Normal"

$r = $ws.Range("E12")
$r.NumberFormat = "@"
$r.Value = "  -0.75%  "
$r.Style = "Normal"

$r = $ws.Range("E13")
$r.NumberFormat = "@"
$r.Value = "  +6.34%  "
$r.Style = "Normal"

$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "2.538.58"
$r.Style = "Normal"

$r = $ws.Range("E14")
$r.NumberFormat = "@"
$r.Value = "  -1.63%  "
$r.Style = "Normal"

$r = $ws.Range("E15")
$r.NumberFormat = "@"
$r.Value = "  -0.98%  "
$r.Style = "Normal"

$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "14.22"
$r.Style = "Normal"

$r = $ws.Range("E16")
$r.NumberFormat = "@"
$r.Value = "  -1.05%  "
$r.Style = "Normal"

$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "42.907.10"
$r.Style = "Normal"

$r = $ws.Range("E17")
$r.NumberFormat = "@"
$r.Value = "  -1.71%  "
$r.Style = "Normal"

$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "0.0₃0991"
$r.Style = "Normal"

$r = $ws.Range("E18")
$r.NumberFormat = "@"
$r.Value = "  +1.51%  "
$r.Style = "Normal"

$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "12.64"
$r.Style = "Normal"

$r = $ws.Range("E19")
$r.NumberFormat = "@"
$r.Value = "  +2.27%  "
$r.Style = "Normal"

$r = $ws.Range("E20")
$r.NumberFormat = "@"
$r.Value = "  -1.17%  "
$r.Style = "Normal"

$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "71.60"
$r.Style = "Normal"

$r = $ws.Range("E21")
$r.NumberFormat = "@"
$r.Value = "  -2.16%  "
$r.Style = "Normal"

$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "253.27"
$r.Style = "Normal"

$r = $ws.Range("E22")
$r.NumberFormat = "@"
$r.Value = "  -4.76%  "
$r.Style = "Normal"

$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "2.94"
$r.Style = "Normal"

$r = $ws.Range("E23")
$r.NumberFormat = "@"
$r.Value = "  +0.19%  "
$r.Style = "Normal"

$r = $ws.Range("E24")
$r.NumberFormat = "@"
$r.Value = "  -4.49%  "
$r.Style = "Normal"

$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "28.74"
$r.Style = "Normal"

$r = $ws.Range("E25")
$r.NumberFormat = "@"
$r.Value = "  -2.21%  "
$r.Style = "Normal"

$r = $ws.Range("E26")
$r.NumberFormat = "@"
$r.Value = "  -0.18%  "
$r.Style = "Normal"

$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "10.26"
$r.Style = "Normal"

$r = $ws.Range("E27")
$r.NumberFormat = "@"
$r.Value = "  +0.13%  "
$r.Style = "Normal"

$r = $ws.Range("E28")
$r.NumberFormat = "@"
$r.Value = "  -1.47%  "
$r.Style = "Normal"

$r = $ws.Range("E29")
$r.NumberFormat = "@"
$r.Value = "  -3.88%  "
$r.Style = "Normal"

$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "6.01"
$r.Style = "Normal"

$r = $ws.Range("E30")
$r.NumberFormat = "@"
$r.Value = "  +0.33%  "
$r.Style = "Normal"

$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "153.82"
$r.Style = "Normal"

$r = $ws.Range("E31")
$r.NumberFormat = "@"
$r.Value = "  +1.09%  "
$r.Style = "Normal"

$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "2.75"
$r.Style = "Normal"

$r = $ws.Range("E32")
$r.NumberFormat = "@"
$r.Value = "  -1.44%  "
$r.Style = "Normal"

$r = $ws.Range("E33")
$r.NumberFormat = "@"
$r.Value = "  -6.80%  "
$r.Style = "Normal"

$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "2.15"
$r.Style = "Normal"

$r = $ws.Range("E34")
$r.NumberFormat = "@"
$r.Value = "  -4.06%  "
$r.Style = "Normal"

$r = $ws.Range("E35")
$r.NumberFormat = "@"
$r.Value = "  -1.20%  "
$r.Style = "Normal"

$r = $ws.Range("E36")
$r.NumberFormat = "@"
$r.Value = "  +6.74%  "
$r.Style = "Normal"

$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "0.114"
$r.Style = "Normal"

$r = $ws.Range("E37")
$r.NumberFormat = "@"
$r.Value = "  -3.04%  "
$r.Style = "Normal"

$r = $ws.Range("E38")
$r.NumberFormat = "@"
$r.Value = "  -0.81%  "
$r.Style = "Normal"

$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "23.16"
$r.Style = "Normal"

$r = $ws.Range("E39")
$r.NumberFormat = "@"
$r.Value = "  -5.17%  "
$r.Style = "Normal"

$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "2.16"
$r.Style = "Normal"

$r = $ws.Range("E40")
$r.NumberFormat = "@"
$r.Value = "  +33.99%  "
$r.Style = "Normal"

$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "3.43"
$r.Style = "Normal"

$r = $ws.Range("E41")
$r.NumberFormat = "@"
$r.Value = "  -2.31%  "
$r.Style = "Normal"

$r = $ws.Range("B42")
$r.NumberFormat = "@"
$r.Value = "RenderToken"
$r.Style = "Normal"

$r = $ws.Range("C42")
$r.NumberFormat = "@"
$r.Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$r.Style = "Normal"

$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "3.89"
$r.Style = "Normal"

$r = $ws.Range("E42")
$r.NumberFormat = "@"
$r.Value = "  +1.01%  "
$r.Style = "Normal"

$r = $ws.Range("B43")
$r.NumberFormat = "@"
$r.Value = "VeChain"
$r.Style = "Normal"

$r = $ws.Range("C43")
$r.NumberFormat = "@"
$r.Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$r.Style = "Normal"

$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "0.0310"
$r.Style = "Normal"

$r = $ws.Range("E43")
$r.NumberFormat = "@"
$r.Value = "  -1.24%  "
$r.Style = "Normal"

$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "2.094.42"
$r.Style = "Normal"

$r = $ws.Range("E44")
$r.NumberFormat = "@"
$r.Value = "  +1.26%  "
$r.Style = "Normal"

$r = $ws.Range("E45")
$r.NumberFormat = "@"
$r.Value = "  +0.14%  "
$r.Style = "Normal"

$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "9.29"
$r.Style = "Normal"

$r = $ws.Range("E46")
$r.NumberFormat = "@"
$r.Value = "  +0.93%  "
$r.Style = "Normal"

$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "85.20"
$r.Style = "Normal"

$r = $ws.Range("E47")
$r.NumberFormat = "@"
$r.Value = "  -3.33%  "
$r.Style = "Normal"

$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "75.88"
$r.Style = "Normal"

$r = $ws.Range("E48")
$r.NumberFormat = "@"
$r.Value = "  +9.63%  "
$r.Style = "Normal"

$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "106.59"
$r.Style = "Normal"

$r = $ws.Range("E49")
$r.NumberFormat = "@"
$r.Value = "  +0.59%  "
$r.Style = "Normal"

$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "2.806.69"
$r.Style = "Normal"

$r = $ws.Range("E50")
$r.NumberFormat = "@"
$r.Value = "  -0.82%  "
$r.Style = "Normal"

$r = $ws.Range("B51")
$r.NumberFormat = "@"
$r.Value = "Stacks"
$r.Style = "Normal"

$r = $ws.Range("C51")
$r.NumberFormat = "@"
$r.Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$r.Style = "Normal"

$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "1.68"
$r.Style = "Normal"

$r = $ws.Range("E51")
$r.NumberFormat = "@"
$r.Value = "  +1.48%  "
$r.Style = "Normal"
